$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old single-letter group codes from column B (rows 3-5: first/mid/last group)
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()

# Add explanatory note in new column C for the group rows
$note = "; blank if there is only one group"
$ws.Range("C3").Value = $note
$ws.Range("C4").Value = $note
$ws.Range("C5").Value = $note

# Size the new column to fit its content
$ws.Columns.Item(3).ColumnWidth = 26.1640625

# Update the selection to the newly added range
$ws.Range("C3:C5").Select()
